$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.111.09'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.651.85'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '218.52'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.5217'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.2646'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '0.06343'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').Value = '4.628'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('D13').Value = '1.649.08'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '1.881.22'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').Value = '0.5598'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').Value = '0.0₅8167'
$ws.Range('D17').Value = '65.44'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '26.114.45'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '4.635'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '10.49'
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('D22').Value = '191.53'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('D23').Value = '5.937'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').Value = '145.53'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').Value = '0.1193'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = '7.230'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('D29').Value = '1.511'
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('D30').Value = '0.05477'
$ws.Range('E30').Value = '  -3.84%  '
$ws.Range('D31').Value = '1.270'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '3.452'
$ws.Range('D33').Value = '3.368'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').Value = '1.560'
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('D35').Value = '0.9521'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').Value = '2.399'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('D38').Value = '0.5642'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = '0.01579'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').Value = '5.865'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').Value = '0.8333'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').Value = '1.029.61'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('D44').Value = '101.16'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').Value = '1.791.63'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('D47').Value = '0.0₈107'
$ws.Range('E47').Value = '  +3.35%  '
$ws.Range('D48').Value = '0.9998'
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '8.006'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.4341'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').Value = '0.05184'
$ws.Range('E51').Value = '  -3.47%  '
